$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grupos")

# Update the stats for "Santa Cruz" (row 13)
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1

# Update the stats for "Tropinha" (row 14)
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = -3
$ws.Range("G14").Value = 1

# Update the active selection to G14
$ws.Range("G14").Select()
